# Generate Report for Handback
#
# - Updates the "Status" text everywhere it appears (Overview!E2:F3 and the
#   "Status" column on the zh-cn / de-de detail sheets) from
#   "Ready for handoff" to "Handed back: in sync with en-US".
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns on the zh-cn and de-de sheets, and
#   hyperlinks the new "Latest Target File" cell the same way column A is
#   hyperlinked.
# - Widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Status text: Overview sheet (zh-cn / de-de status columns) ----
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# ---- Status text: per-language "Status" column ----
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---- zh-cn sheet: fill in target/handback info for both rows ----
$wsZhCn.Range("J2").Value = "6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.a9ddb40dba1d75e53eee9622e49518ae5df79f5a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-23 10:29:54"
$wsZhCn.Range("J3").Value = "e2330687-05fe-4915-9e0d-e4a06ac0015a.d0c26e6a8420be42ef3cddb8121ca0d740500521.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-23 10:29:54"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/709171ac856c3bbdce6cd936cba98782f7d4c373/e2e/6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.md", "", "", "6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/709171ac856c3bbdce6cd936cba98782f7d4c373/e2e/e2330687-05fe-4915-9e0d-e4a06ac0015a.md", "", "", "e2330687-05fe-4915-9e0d-e4a06ac0015a.md")

# ---- de-de sheet: fill in target/handback info for both rows ----
$wsDeDe.Range("J2").Value = "6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.a9ddb40dba1d75e53eee9622e49518ae5df79f5a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-23 10:30:11"
$wsDeDe.Range("J3").Value = "e2330687-05fe-4915-9e0d-e4a06ac0015a.d0c26e6a8420be42ef3cddb8121ca0d740500521.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-23 10:30:11"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/709171ac856c3bbdce6cd936cba98782f7d4c373/e2e/6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.md", "", "", "6de7e0e2-1eca-4f1f-a0ff-5a100bea1eba.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/709171ac856c3bbdce6cd936cba98782f7d4c373/e2e/e2330687-05fe-4915-9e0d-e4a06ac0015a.md", "", "", "e2330687-05fe-4915-9e0d-e4a06ac0015a.md")

# ---- Column widths (values now need more room) ----
# Overview: zh-cn / de-de status columns
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status column
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668

# zh-cn / de-de: Latest Target File / Latest Handback File columns
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
